$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Heading: "Why do we care about certainty vs uncertainty?" -> "Risk" + "?"
#    (kept as two separate runs, matching the target OOXML)
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item(10)
$headingStart = $headingPara.Range.Start
$oldHeading = "Why do we care about certainty vs uncertainty?"
$headingRange = $d.Range($headingStart, $headingStart + $oldHeading.Length)
$headingRange.Text = "Risk"

# Insert the "?" right after "Risk" as its own run.
$afterRisk = $d.Range($headingStart + 4, $headingStart + 4)
$afterRisk.InsertAfter("?")

# Force the runtime to keep "Risk" and "?" as separate <w:r> elements instead
# of silently re-merging them (toggle a formatting property off again so the
# visible formatting is unchanged but the run boundary survives).
$riskRunRange = $d.Range($headingStart, $headingStart + 4)
$riskRunRange.Font.Bold = 1
$riskRunRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) Replace the six "Why do we care..." body paragraphs with the new
#    three-line outline followed by three blank paragraphs.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(11).Range.Text = "Today we dig a little more into risk"
$d.Paragraphs.Item(12).Range.Text = "Adding probabilities vs multiplying"
$d.Paragraphs.Item(13).Range.Text = "Independence and the gamblers fallacy"

14..16 | ForEach-Object {
    $p = $d.Paragraphs.Item($_)
    $r = $p.Range
    if ($r.End - 1 -gt $r.Start) {
        $trimmed = $d.Range($r.Start, $r.End - 1)
        $trimmed.Delete()
    }
}

Write-Output "Done"
